$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, matching the style of the existing
# header cells (E1, etc.) by copying formats only (xlPasteFormats = -4122)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Per-row "time_taken" timestamps (text), recorded to microsecond precision
$timestamps = @(
    "2021-10-05 10:50:25.822205",
    "2021-10-05 10:50:25.822217",
    "2021-10-05 10:50:25.822220",
    "2021-10-05 10:50:25.822223",
    "2021-10-05 10:50:25.822226",
    "2021-10-05 10:50:25.822228",
    "2021-10-05 10:50:25.822231",
    "2021-10-05 10:50:25.822233",
    "2021-10-05 10:50:25.822236",
    "2021-10-05 10:50:25.822239",
    "2021-10-05 10:50:25.822241",
    "2021-10-05 10:50:25.822244",
    "2021-10-05 10:50:25.822246",
    "2021-10-05 10:50:25.822249",
    "2021-10-05 10:50:25.822251",
    "2021-10-05 10:50:25.822254",
    "2021-10-05 10:50:25.822260",
    "2021-10-05 10:50:25.822263",
    "2021-10-05 10:50:25.822266",
    "2021-10-05 10:50:25.822268",
    "2021-10-05 10:50:25.822271",
    "2021-10-05 10:50:25.822273",
    "2021-10-05 10:50:25.822275",
    "2021-10-05 10:50:25.822278",
    "2021-10-05 10:50:25.822281",
    "2021-10-05 10:50:25.822283",
    "2021-10-05 10:50:25.822286",
    "2021-10-05 10:50:25.822288",
    "2021-10-05 10:50:25.822291",
    "2021-10-05 10:50:25.822293",
    "2021-10-05 10:50:25.822296",
    "2021-10-05 10:50:25.822298",
    "2021-10-05 10:50:25.822301",
    "2021-10-05 10:50:25.822304",
    "2021-10-05 10:50:25.822307",
    "2021-10-05 10:50:25.822309",
    "2021-10-05 10:50:25.822312",
    "2021-10-05 10:50:25.822315",
    "2021-10-05 10:50:25.822317",
    "2021-10-05 10:50:25.822320",
    "2021-10-05 10:50:25.822323",
    "2021-10-05 10:50:25.822326",
    "2021-10-05 10:50:25.822328",
    "2021-10-05 10:50:25.822331",
    "2021-10-05 10:50:25.822333",
    "2021-10-05 10:50:25.822335",
    "2021-10-05 10:50:25.822338",
    "2021-10-05 10:50:25.822340",
    "2021-10-05 10:50:25.822343",
    "2021-10-05 10:50:25.822346",
    "2021-10-05 10:50:25.822348",
    "2021-10-05 10:50:25.822350",
    "2021-10-05 10:50:25.822353",
    "2021-10-05 10:50:25.822356",
    "2021-10-05 10:50:25.822358",
    "2021-10-05 10:50:25.822361",
    "2021-10-05 10:50:25.822363",
    "2021-10-05 10:50:25.822366",
    "2021-10-05 10:50:25.822368",
    "2021-10-05 10:50:25.822371",
    "2021-10-05 10:50:25.822373",
    "2021-10-05 10:50:25.822376",
    "2021-10-05 10:50:25.822378",
    "2021-10-05 10:50:25.822381",
    "2021-10-05 10:50:25.822384",
    "2021-10-05 10:50:25.822387",
    "2021-10-05 10:50:25.822389",
    "2021-10-05 10:50:25.822392",
    "2021-10-05 10:50:25.822394",
    "2021-10-05 10:50:25.822397",
    "2021-10-05 10:50:25.822399",
    "2021-10-05 10:50:25.822402",
    "2021-10-05 10:50:25.822404",
    "2021-10-05 10:50:25.822407",
    "2021-10-05 10:50:25.822409",
    "2021-10-05 10:50:25.822412",
    "2021-10-05 10:50:25.822417",
    "2021-10-05 10:50:25.822420",
    "2021-10-05 10:50:25.822422",
    "2021-10-05 10:50:25.822425",
    "2021-10-05 10:50:25.822428",
    "2021-10-05 10:50:25.822430",
    "2021-10-05 10:50:25.822432",
    "2021-10-05 10:50:25.822435",
    "2021-10-05 10:50:25.822438",
    "2021-10-05 10:50:25.822440",
    "2021-10-05 10:50:25.822443",
    "2021-10-05 10:50:25.822445",
    "2021-10-05 10:50:25.822448",
    "2021-10-05 10:50:25.822450",
    "2021-10-05 10:50:25.822453",
    "2021-10-05 10:50:25.822455",
    "2021-10-05 10:50:25.822459",
    "2021-10-05 10:50:25.822462",
    "2021-10-05 10:50:25.822465",
    "2021-10-05 10:50:25.822467",
    "2021-10-05 10:50:25.822470",
    "2021-10-05 10:50:25.822473",
    "2021-10-05 10:50:25.822475",
    "2021-10-05 10:50:25.822478",
    "2021-10-05 10:50:25.822480",
    "2021-10-05 10:50:25.822483",
    "2021-10-05 10:50:25.822485",
    "2021-10-05 10:50:25.822488",
    "2021-10-05 10:50:25.822491",
    "2021-10-05 10:50:25.822493",
    "2021-10-05 10:50:25.822496",
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
